$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item("Proportional - Constant")

# --- Add the new "Proportional - Proportional" sheet as the last tab.
#     Duplicate "Proportional - Constant" (right after it, while it is still
#     in its original/template state) so the new sheet inherits the exact
#     same layout / column widths / number formats, then rename it and blank
#     its data back out to the simulation template (all zeros except the
#     first edge-count result that has been entered so far). ---
$ws4.Copy($null, $ws4)
$ws5 = $wb.Worksheets.Item("Proportional - Constant (2)")
$ws5.Name = "Proportional - Proportional"

$ws5.Range("C2:E11").Value = 0
$ws5.Range("C2").Value = 59890

$ws5.Range("C3").Select()

# --- Finish the "Proportional - Constant" sheet: the last three simulation
#     rows (9, 10, 11) were still zero-filled placeholders; fill in the real
#     results so the averages in row 12 recompute. ---
$ws4.Range("D9").Value = 28806
$ws4.Range("E9").Value = 28869

$ws4.Range("C10").Value = 59877
$ws4.Range("D10").Value = 28802
$ws4.Range("E10").Value = 28869

$ws4.Range("C11").Value = 59889
$ws4.Range("D11").Value = 28804
$ws4.Range("E11").Value = 28871

# Sheet4 is done being edited: selection moves from the last-edited cell to
# the whole data block, and it stops being the tab-selected sheet (the new
# sheet5 is now the active/selected tab).
$ws4.Range("A1:E12").Select()
$ws5.Activate()
